# Updated cryptos list on Fri May 26 07:09:20 UTC 2023 with GitHub Actions
#
# The sheet lists coins in columns B (Coin), C (Link), D (Price), E (Volume 1h).
# All D/E values (and, for 3 rows, B/C too) are stored as plain TEXT strings
# (not numbers) even when the text happens to look numeric (e.g. "1.009").
# Writing a numeric-looking string straight into .Value would make Excel
# auto-coerce it into a real number (and reformat it, e.g. "1.009" ->
# 1.0089999999999999), which would not match the source data. To keep such
# values as literal text we prefix them with a leading apostrophe (the
# classic "force text" entry trick) and then reset the cell Style back to
# "Normal" so no stray number-format/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Addr,
        [string]$NewValue
    )

    # Plain decimal/integer (with optional sign) -> Excel would parse this
    # as a real number if assigned directly.
    $looksNumeric = $NewValue -match '^[+-]?\d+(\.\d+)?$'

    $cell = $ws.Range($Addr)
    if ($looksNumeric) {
        # Force text entry so Excel doesn't coerce this into a real number.
        $cell.Value = "'" + $NewValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $NewValue
    }
}

# --- Row 2: Bitcoin ---
Set-TextCell "D2" "26.598.91"
Set-TextCell "E2" "  +0.71%  "

# --- Row 3: Ethereum ---
Set-TextCell "D3" "1.822.73"
Set-TextCell "E3" "  +1.47%  "

# --- Row 4: TetherUSD ---
Set-TextCell "D4" "1.009"
Set-TextCell "E4" "  +0.06%  "

# --- Row 5: USDC ---
Set-TextCell "D5" "1.007"
Set-TextCell "E5" "  -0.03%  "

# --- Row 6: BNB ---
Set-TextCell "D6" "307.67"
Set-TextCell "E6" "  -0.18%  "

# --- Row 7: XRP ---
Set-TextCell "D7" "0.4647"
Set-TextCell "E7" "  +2.57%  "

# --- Row 8: Cardano ---
Set-TextCell "D8" "0.3602"
Set-TextCell "E8" "  +0.20%  "

# --- Row 9: Dogecoin ---
Set-TextCell "D9" "0.07128"
Set-TextCell "E9" "  +0.07%  "

# --- Row 10: Polygon ---
Set-TextCell "D10" "0.9009"
Set-TextCell "E10" "  +1.54%  "

# --- Row 11: TRON ---
Set-TextCell "D11" "0.07770"
Set-TextCell "E11" "  -0.76%  "

# --- Row 12: Solana ---
Set-TextCell "D12" "19.38"
Set-TextCell "E12" "  -0.73%  "

# --- Row 13: WrappedEther ---
Set-TextCell "D13" "1.812.47"
Set-TextCell "E13" "  +0.62%  "

# --- Row 14: Polkadot ---
Set-TextCell "D14" "5.260"
Set-TextCell "E14" "  -0.41%  "

# --- Row 15: Chainlink ---
Set-TextCell "D15" "6.320"
Set-TextCell "E15" "  -0.27%  "

# --- Row 16: Litecoin ---
Set-TextCell "D16" "87.34"
Set-TextCell "E16" "  +2.78%  "

# --- Row 17: BinanceUSD ---
Set-TextCell "D17" "1.009"
Set-TextCell "E17" "  -0.09%  "

# --- Row 18: ShibaInu ---
Set-TextCell "D18" "0.000008553"
Set-TextCell "E18" "  -0.34%  "

# --- Row 19: Dai ---
Set-TextCell "E19" "  -0.05%  "

# --- Row 20: WrappedBTC ---
Set-TextCell "D20" "26.648.32"
Set-TextCell "E20" "  +0.76%  "

# --- Row 21: Avalanche ---
Set-TextCell "D21" "14.15"
Set-TextCell "E21" "  -1.02%  "

# --- Row 22: Uniswap ---
Set-TextCell "D22" "5.014"
Set-TextCell "E22" "  +0.43%  "

# --- Row 23: Cosmos ---
Set-TextCell "E23" "  -0.02%  "

# --- Row 24: Toncoin ---
Set-TextCell "D24" "1.919"
Set-TextCell "E24" "  -3.38%  "

# --- Row 25: Monero ---
Set-TextCell "D25" "152.08"
Set-TextCell "E25" "  -0.37%  "

# --- Row 26: EthereumClassic ---
Set-TextCell "D26" "17.90"
Set-TextCell "E26" "  -0.24%  "

# --- Row 27: LidoDAOToken ---
Set-TextCell "D27" "1.981"
Set-TextCell "E27" "  -3.35%  "

# --- Row 28: BitcoinCash ---
Set-TextCell "D28" "113.76"
Set-TextCell "E28" "  +1.48%  "

# --- Row 29: InternetComputer(DFINITY) ---
Set-TextCell "D29" "4.802"
Set-TextCell "E29" "  -1.63%  "

# --- Row 30: Stellar ---
Set-TextCell "D30" "0.08808"
Set-TextCell "E30" "  +1.69%  "

# --- Row 31: HuobiToken ---
Set-TextCell "D31" "3.137"
Set-TextCell "E31" "  +2.92%  "

# --- Row 32: ImmutableX ---
Set-TextCell "D32" "0.7313"
Set-TextCell "E32" "  +0.35%  "

# --- Row 33: was RenderToken -> now Filecoin ---
Set-TextCell "B33" "Filecoin"
Set-TextCell "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D33" "4.435"
Set-TextCell "E33" "  -0.64%  "

# --- Row 34: was Filecoin -> now ARBITRUM ---
Set-TextCell "B34" "ARBITRUM"
Set-TextCell "C34" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D34" "1.137"
Set-TextCell "E34" "  +2.13%  "

# --- Row 35: was ARBITRUM -> now RenderToken ---
Set-TextCell "B35" "RenderToken"
Set-TextCell "C35" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D35" "2.714"
Set-TextCell "E35" "  -1.66%  "

# --- Row 36: TrustWalletToken ---
Set-TextCell "D36" "1.072"
Set-TextCell "E36" "  -0.36%  "

# --- Row 37: VeChain ---
Set-TextCell "D37" "0.01921"
Set-TextCell "E37" "  -1.06%  "

# --- Row 38: MXToken ---
Set-TextCell "E38" "  +1.29%  "

# --- Row 39: Hedera ---
Set-TextCell "D39" "0.05116"
Set-TextCell "E39" "  +0.10%  "

# --- Row 40: FraxShare ---
Set-TextCell "D40" "6.895"
Set-TextCell "E40" "  -0.19%  "

# --- Row 41: TheSandbox ---
Set-TextCell "D41" "0.5050"
Set-TextCell "E41" "  -2.37%  "

# --- Row 42: Algorand ---
Set-TextCell "D42" "0.1494"
Set-TextCell "E42" "  -1.83%  "

# --- Row 43: Aptos ---
Set-TextCell "D43" "7.987"
Set-TextCell "E43" "  -0.38%  "

# --- Row 44: PaxDollar ---
Set-TextCell "E44" "  -0.04%  "

# --- Row 45: Decentraland ---
Set-TextCell "D45" "0.4649"
Set-TextCell "E45" "  -0.59%  "

# --- Row 46: EnergySwap ---
Set-TextCell "D46" "9.960"
Set-TextCell "E46" "  +0.85%  "

# --- Row 47: Quant ---
Set-TextCell "D47" "98.11"
Set-TextCell "E47" "  -2.48%  "

# --- Row 48: NEARProtocol ---
Set-TextCell "D48" "1.557"
Set-TextCell "E48" "  -1.93%  "

# --- Row 49: Cronos ---
Set-TextCell "D49" "0.05982"
Set-TextCell "E49" "  +0.03%  "

# --- Row 50: Aave ---
Set-TextCell "D50" "63.67"
Set-TextCell "E50" "  -1.29%  "

# --- Row 51: Elrond ---
Set-TextCell "D51" "35.80"
Set-TextCell "E51" "  -0.97%  "
